$wb = $excel.ActiveWorkbook

# Update the "Code Quality" sheet values (Sonarqube data changed)
$ws = $wb.Worksheets.Item("Code Quality")
$ws.Range("C4").Value = 80
$ws.Range("C5").Value = 72
$ws.Range("D5").Value = 65
$ws.Range("C6").Value = 2.5
$ws.Range("C7").Value = 3

# Move selection / active sheet to "Code Quality" (was "Physics")
$ws.Activate()
$ws.Range("G5").Select()
